# Swap the match-data (columns B:AC) between specific pairs of rows.
# Column A (the running index/id) stays attached to its own row; only the
# match record fields (B..AC) move between the two rows in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(42, 43),
    @(55, 56),
    @(61, 62),
    @(65, 66),
    @(74, 75),
    @(88, 89),
    @(100, 101),
    @(119, 120),
    @(124, 125),
    @(155, 156),
    @(170, 171),
    @(182, 183),
    @(200, 201),
    @(241, 242),
    @(245, 246),
    @(262, 263),
    @(271, 272)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 29))
    $range2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 29))

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
